$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.196.18'
$ws.Range('E2').Value = '  +1.68%  '

$ws.Range('D3').Value = '3.924.63'
$ws.Range('E3').Value = '  +0.50%  '

$ws.Range('E4').Value = '  +0.08%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '486.46'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.71%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.88'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.69%  '

$ws.Range('E7').Value = '  -0.25%  '

$ws.Range('E8').Value = '  -0.13%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.728'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.21%  '

$ws.Range('E10').Value = '  +3.48%  '

$ws.Range('E11').Value = '  +6.12%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '42.59'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.84%  '

$ws.Range('E13').Value = '  +1.98%  '

$ws.Range('D14').Value = '4.549.93'
$ws.Range('E14').Value = '  +0.51%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.84'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.93%  '

$ws.Range('D16').Value = '3.938.21'
$ws.Range('E16').Value = '  +0.70%  '

$ws.Range('E17').Value = '  -0.19%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.91'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.44%  '

$ws.Range('E19').Value = '  -1.66%  '

$ws.Range('D20').Value = '68.312.43'
$ws.Range('E20').Value = '  +1.49%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '446.28'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.38%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.85'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.00%  '

$ws.Range('E23').Value = '  +1.49%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '88.74'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.10%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.42'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +13.02%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.85'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +12.87%  '

$ws.Range('E27').Value = '  +2.89%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '38.77'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.51%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.88'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.00%  '

$ws.Range('B30').Value = 'Cosmos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '13.45'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.33%  '

$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '691.10'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.54%  '

$ws.Range('E32').Value = '  -0.56%  '

$ws.Range('E33').Value = '  +5.11%  '

$ws.Range('D34').Value = '0.0₃0948'
$ws.Range('E34').Value = '  +18.99%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '41.67'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.10%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '59.07'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.62%  '

$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.70'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.83%  '

$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.150'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.23%  '

$ws.Range('E39').Value = '  -0.05%  '

$ws.Range('E40').Value = '  +0.74%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.88'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +15.77%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.13'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.18%  '

$ws.Range('E43').Value = '  +10.39%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.94'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.84%  '

$ws.Range('E45').Value = '  +1.63%  '

$ws.Range('E46').Value = '  -0.09%  '

$ws.Range('E47').Value = '  +0.29%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.15'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.06%  '

$ws.Range('E49').Value = '  +1.58%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '145.27'
$ws.Range('D50').Style = 'Normal'

$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0342'
$ws.Range('E51').Value = '  +37.59%  '
